# heuristic_results.xlsx update:
#  - Summary sheet: new Total Cost / Air Shipping Cost values, plus new rows for
#    Total Holding Cost, Total Fixed Cost, Total Quantity Ordered, Total Ocean
#    Shipping Quantity, Total Air Shipping Quantity (Number of Products / Number
#    of Periods shift down to make room).
#  - A brand-new "Quantities Summary" sheet is inserted right after "Summary"
#    (it reuses the old "Purchasing Costs" sheet position/content, restructured
#    with Ocean/Air quantity + volume columns).
#  - A brand-new "Purchasing Costs" sheet is inserted after "Quantities Summary"
#    (holding the unit-cost / quantity / purchasing+holding+fixed cost table).
#  - Orders sheet is untouched.
#  - Inventory sheet gets a handful of ending-inventory values updated.
#  - Containers sheet is untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("B2").Value = 16515758
$summary.Range("B4").Value = 8418

# carry plain data-row styling down onto the newly-created rows (A5:B5 is an
# unstyled data row, same as the rest of the metric/value pairs) before
# filling in their text/values
$summary.Range("A5:B5").Copy($summary.Range("A6:B6"))
$summary.Range("A5:B5").Copy($summary.Range("A7:B7"))
$summary.Range("A5:B5").Copy($summary.Range("A8:B8"))
$summary.Range("A5:B5").Copy($summary.Range("A9:B9"))
$summary.Range("A5:B5").Copy($summary.Range("A10:B10"))
$summary.Range("A5:B5").Copy($summary.Range("A11:B11"))
$summary.Range("A5:B5").Copy($summary.Range("A12:B12"))

$summary.Range("A6").Value = "Total Holding Cost"
$summary.Range("B6").Value = 1256400

$summary.Range("A7").Value = "Total Fixed Cost"
$summary.Range("B7").Value = 2440

$summary.Range("A8").Value = "Number of Products"
$summary.Range("B8").Value = 10

$summary.Range("A9").Value = "Number of Periods"
$summary.Range("B9").Value = 6

$summary.Range("A10").Value = "Total Quantity Ordered (All Products)"
$summary.Range("B10").Value = 2153

$summary.Range("A11").Value = "Total Ocean Shipping Quantity"
$summary.Range("B11").Value = 1967

$summary.Range("A12").Value = "Total Air Shipping Quantity"
$summary.Range("B12").Value = 186

Write-Host "Summary sheet updated"

# ---------------------------------------------------------------------------
# 2. "Purchasing Costs" -> "Quantities Summary" (restructured columns)
# ---------------------------------------------------------------------------
# Old layout: A=Product B=Unit Cost C=Total Quantity D=Total Purchasing Cost
# New layout: A=Product B=Total Quantity C=Ocean Shipping Quantity
#             D=Air Shipping Quantity E=Unit Cost F=Total Purchasing Cost
#             G=Total Volume H=Ocean Volume I=Air Volume
$qs = $wb.Worksheets.Item("Purchasing Costs")
$qs.Name = "Quantities Summary"

# extend header styling (A1:D1 already styled) across the new header cells
$qs.Range("D1").Copy($qs.Range("E1:I1"))

$qs.Range("A1").Value = "Product"
$qs.Range("B1").Value = "Total Quantity"
$qs.Range("C1").Value = "Ocean Shipping Quantity"
$qs.Range("D1").Value = "Air Shipping Quantity"
$qs.Range("E1").Value = "Unit Cost"
$qs.Range("F1").Value = "Total Purchasing Cost"
$qs.Range("G1").Value = "Total Volume"
$qs.Range("H1").Value = "Ocean Volume"
$qs.Range("I1").Value = "Air Volume"

$qsData = @(
    @(1, 38, 38, 0, 5000, 190000, 2.774, 2.774, 0),
    @(2, 45, 45, 0, 2000, 90000, 0.225, 0.225, 0),
    @(3, 282, 282, 0, 9000, 2538000, 12.126, 12.126, 0),
    @(4, 252, 252, 0, 9000, 2268000, 15.876, 15.876, 0),
    @(5, 116, 116, 0, 2000, 232000, 5.22, 5.22, 0),
    @(6, 154, 154, 0, 9000, 1386000, 13.244, 13.244, 0),
    @(7, 172, 172, 0, 7000, 1204000, 13.588, 13.588, 0),
    @(8, 261, 261, 0, 5000, 1305000, 21.402, 21.402, 0),
    @(9, 94, 94, 0, 9000, 846000, 6.392, 6.392, 0),
    @(10, 739, 553, 186, 7000, 5173000, 72.422, 54.194, 18.228)
)

$row = 2
foreach ($rec in $qsData) {
    for ($col = 1; $col -le 9; $col++) {
        $qs.Cells.Item($row, $col).Value = $rec[$col - 1]
    }
    $row = $row + 1
}

Write-Host "Quantities Summary sheet built"

# ---------------------------------------------------------------------------
# 3. New "Purchasing Costs" sheet (inserted right after "Quantities Summary",
#    i.e. right before "Orders")
# ---------------------------------------------------------------------------
$orders = $wb.Worksheets.Item("Orders")
$pc = $wb.Worksheets.Add($orders)
$pc.Name = "Purchasing Costs"

# reuse the bold/bordered/centered header style from the Quantities Summary
# sheet's header row so the new sheet matches the workbook's look
$qs.Range("A1:H1").Copy($pc.Range("A1:H1"))

$pc.Range("A1").Value = "Product"
$pc.Range("B1").Value = "Unit Cost"
$pc.Range("C1").Value = "Total Quantity"
$pc.Range("D1").Value = "Ocean Quantity"
$pc.Range("E1").Value = "Air Quantity"
$pc.Range("F1").Value = "Total Purchasing Cost"
$pc.Range("G1").Value = "Total Holding Cost"
$pc.Range("H1").Value = "Total Fixed Cost"

$pcData = @(
    @(1, 5000, 38, 38, 0, 190000, 125640, 244),
    @(2, 2000, 45, 45, 0, 90000, 125640, 244),
    @(3, 9000, 282, 282, 0, 2538000, 125640, 244),
    @(4, 9000, 252, 252, 0, 2268000, 125640, 244),
    @(5, 2000, 116, 116, 0, 232000, 125640, 244),
    @(6, 9000, 154, 154, 0, 1386000, 125640, 244),
    @(7, 7000, 172, 172, 0, 1204000, 125640, 244),
    @(8, 5000, 261, 261, 0, 1305000, 125640, 244),
    @(9, 9000, 94, 94, 0, 846000, 125640, 244),
    @(10, 7000, 739, 553, 186, 5173000, 125640, 244)
)

$row = 2
foreach ($rec in $pcData) {
    for ($col = 1; $col -le 8; $col++) {
        $pc.Cells.Item($row, $col).Value = $rec[$col - 1]
    }
    $row = $row + 1
}

Write-Host "Purchasing Costs sheet built"

# ---------------------------------------------------------------------------
# 4. Inventory sheet: a handful of ending-inventory values change (period 2
#    and period 3 rows for products 2, 3, 4, 6, 7, 9)
# ---------------------------------------------------------------------------
$inv = $wb.Worksheets.Item("Inventory")

$inv.Range("C13").Value = 357
$inv.Range("C15").Value = 258
$inv.Range("C17").Value = 356
$inv.Range("C18").Value = 212
$inv.Range("C20").Value = 489
$inv.Range("C24").Value = 166
$inv.Range("C27").Value = 272
$inv.Range("C28").Value = 238
$inv.Range("C30").Value = 335

Write-Host "Inventory sheet updated"

# ---------------------------------------------------------------------------
# 5. Orders and Containers sheets are left untouched.
# ---------------------------------------------------------------------------

Write-Host "All sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
